$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as literal text
# (e.g. "1.00", "0.692"). Excel auto-converts such strings to numbers on
# assignment, so pre-format the affected cells as Text to keep them as
# strings, matching the source data exactly.
$textCells = @(
    'D4', 'D6', 'D7', 'D8', 'D9', 'D11', 'D13', 'D14', 'D16', 'D19', 'D21', 'D22', 'D23', 'D24', 'D26', 'D27', 'D29', 'D31', 'D33', 'D38', 'D42', 'D43', 'D44', 'D47', 'D48', 'D50', 'D51'
)
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '35.884.69'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '1.891.81'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').Value = '0.692'
$ws.Range('E6').Value = '  +0.74%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '43.23'
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('D9').Value = '56.88'
$ws.Range('E9').Value = '  +9.46%  '
$ws.Range('E10').Value = '  +1.89%  '
$ws.Range('D11').Value = '0.0755'
$ws.Range('E11').Value = '  +2.58%  '
$ws.Range('E12').Value = '  +1.46%  '
$ws.Range('D13').Value = '14.94'
$ws.Range('E13').Value = '  +14.17%  '
$ws.Range('D14').Value = '0.793'
$ws.Range('E14').Value = '  +8.07%  '
$ws.Range('D15').Value = '2.164.84'
$ws.Range('E15').Value = '  +0.30%  '
$ws.Range('D16').Value = '5.07'
$ws.Range('E16').Value = '  +2.43%  '
$ws.Range('D17').Value = '1.888.03'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '35.830.66'
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('D19').Value = '73.66'
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('D20').Value = '0.0₃0832'
$ws.Range('E20').Value = '  +1.45%  '
$ws.Range('D21').Value = '247.31'
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('D22').Value = '13.12'
$ws.Range('E22').Value = '  +2.54%  '
$ws.Range('D23').Value = '5.18'
$ws.Range('E23').Value = '  +5.01%  '
$ws.Range('D24').Value = '2.68'
$ws.Range('E24').Value = '  +4.75%  '
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').Value = '2.15'
$ws.Range('E26').Value = '  -1.05%  '
$ws.Range('D27').Value = '166.31'
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('E28').Value = '  +2.71%  '
$ws.Range('D29').Value = '18.42'
$ws.Range('E29').Value = '  +0.87%  '
$ws.Range('E30').Value = '  +0.70%  '
$ws.Range('D31').Value = '4.43'
$ws.Range('E31').Value = '  +4.76%  '
$ws.Range('E32').Value = '  +5.00%  '
$ws.Range('D33').Value = '4.28'
$ws.Range('E33').Value = '  +1.99%  '
$ws.Range('E34').Value = '  +17.70%  '
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('E36').Value = '  -14.39%  '
$ws.Range('E37').Value = '  +0.88%  '
$ws.Range('D38').Value = '0.0765'
$ws.Range('E38').Value = '  +10.73%  '
$ws.Range('E39').Value = '  -2.20%  '
$ws.Range('E40').Value = '  +7.05%  '
$ws.Range('E41').Value = '  +2.07%  '
$ws.Range('D42').Value = '17.00'
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('D43').Value = '1.08'
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').Value = '14.08'
$ws.Range('E44').Value = '  +16.55%  '
$ws.Range('D45').Value = '1.313.17'
$ws.Range('E45').Value = '  +2.03%  '
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('D47').Value = '0.0810'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').Value = '2.41'
$ws.Range('E48').Value = '  +0.26%  '
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('D50').Value = '6.33'
$ws.Range('E50').Value = '  +0.88%  '
$ws.Range('D51').Value = '42.59'
$ws.Range('E51').Value = '  -1.23%  '
